{"js": "// Applies the \"6000 -> 5000\", \"20000 -> 10000\", and the re-run test-loss /\n// test-accuracy figures (0.3801 -> 0.7565, 0.9457 -> 0.8598) edits described\n// by the diff. The numbers changed consistently everywhere they appear\n// (the table's \"Test Set\"/\"Training Set\" columns and the narrative\n// paragraphs that restate the same sample counts and result figures), so a\n// handful of scoped, literal search & replace passes reproduce the edit.\n\nconst body = context.document.body;\n\nasync function replaceAll(findText, replacement) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Training/Test-set sample counts: every \"20000\" becomes \"10000\" and every\n// \"6000\" becomes \"5000\" (table cells + the narrative mentions of the same\n// numbers).\nawait replaceAll(\"20000\", \"10000\");\nawait replaceAll(\"6000\", \"5000\");\n\n// Re-run results: the final model now reports a higher loss / lower\n// accuracy. Handle the narrative sentence's \"test loss:0.3801\" -> \"test\n// loss of 0.7565\" wording tweak first (it also drops the colon in favor of\n// \"of \"), then sweep the remaining occurrences of the old figures (the\n// summary table cell \"loss: 0.3801 - acc: 0.9457\").\nawait replaceAll(\"test loss:0.3801\", \"test loss of 0.7565\");\nawait replaceAll(\"0.3801\", \"0.7565\");\nawait replaceAll(\"0.9457\", \"0.8598\");\n", "ps1": "# Applies the \"6000 -> 5000\", \"20000 -> 10000\", and the re-run test-loss /\n# test-accuracy figures (0.3801 -> 0.7565, 0.9457 -> 0.8598) edits described\n# by the diff. The numbers changed consistently everywhere they appear (the\n# table's \"Test Set\"/\"Training Set\" columns and the narrative paragraphs\n# that restate the same sample counts and result figures), so a handful of\n# scoped, literal Find/Replace-All passes over the whole document\n# reproduce the edit.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replText, 2) | Out-Null\n}\n\n# Training/Test-set sample counts: every \"20000\" becomes \"10000\" and every\n# \"6000\" becomes \"5000\" (table cells + the narrative mentions of the same\n# numbers).\nReplace-AllText \"20000\" \"10000\"\nReplace-AllText \"6000\" \"5000\"\n\n# Re-run results: the final model now reports a higher loss / lower\n# accuracy. Handle the narrative sentence's \"test loss:0.3801\" -> \"test\n# loss of 0.7565\" wording tweak first (it also drops the colon in favor of\n# \"of \"), then sweep the remaining occurrences of the old figures (the\n# summary table cell \"loss: 0.3801 - acc: 0.9457\").\nReplace-AllText \"test loss:0.3801\" \"test loss of 0.7565\"\nReplace-AllText \"0.3801\" \"0.7565\"\nReplace-AllText \"0.9457\" \"0.8598\"\n"}
